$d = $word.ActiveDocument

# Start from the end of the document's last paragraph and append the new content.
$last = $d.Paragraphs.Last.Range
$last.InsertParagraphAfter()

$heading = $d.Paragraphs.Last.Range
$heading.Style = "Heading 2"
$heading.Text = "ProductViewpoint"
$heading.InsertParagraphAfter()

$p1 = $d.Paragraphs.Last.Range
$p1.Style = "Normal"
$p1.Text = "U product viewpointu možemo vidjeti kako se proizvod kreće iza kulise."
$p1.InsertParagraphAfter()

$p2 = $d.Paragraphs.Last.Range
$p2.Style = "Normal"
$p2.Text = "Za svaku narudžbu na raspolaganju stoji služba za korisnike te logistički odjel koji se brine oko dovoza hrane iz centralne kuhinje, dovoz proizvoda kupljenih od OPG-ova te posljednje i najbitnije dostave same narudžbe klijentu."
$p2.InsertParagraphAfter()

$p3 = $d.Paragraphs.Last.Range
$p3.Style = "Normal"
$p3.Text = "Služba za korisnike nam stoji na raspolaganju za pritužbe, pomoć pri naručivanju proizvoda, probleme se narudžbom te bilo kakve druge upite."
$p3.InsertParagraphAfter()

$p4 = $d.Paragraphs.Last.Range
$p4.Style = "Normal"
$p4.Text = "Te na kraju imamo i poslovnu ulogu mnogobrojnih OPG-ova koji su bili spomenuti u logističkom odjelu, koji prodaju proizvode našoj tvrtci."
$p4.InsertParagraphAfter()

$p5 = $d.Paragraphs.Last.Range
$p5.Style = "Normal"
$p5.Text = "Također imamo i ugovor koji tvrtku veže za dostavljanje narudžbe kupcu."
